$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update arrival count (B2): 38 -> 40
$ws.Range("B2").Value = 40

# Update local minima count (D2): 2 -> 0
$ws.Range("D2").Value = 0

# Update arrival probability (B5): 0.95 -> 1
$ws.Range("B5").Value = 1

# Update minima probability (D5): 0.05 -> 0
$ws.Range("D5").Value = 0
